$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BUS1 column (A3:A9) values per completed Monte Carlo simulation
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 4
$ws.Range("A9").Value = 4

# Move the active selection to A9
$ws.Range("A9").Select()
